$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Screen"
$ws.Range("B1").Value = "Issue"
$ws.Range("A2").Value = "Edit Center"
$ws.Range("B2").Value = """Edit Center"" and ""Add/Edit Center"" both appear on the screen. (Remove ""Edit Center"".)"
$ws.Range("A3").Value = "Edit Center"
$ws.Range("B3").Value = "A blank drop down appears under the City textbox."
$ws.Range("A4").Value = "Manage Volunteers"
$ws.Range("B4").Value = """Manage Volunteers"" appears twice on the screen. (Should appear under breadcrumbs.)"
$ws.Range("A5").Value = "Edit Center"
$ws.Range("B5").Value = "Address field is truncated on save."
$ws.Range("A6").Value = "Edit Task"
$ws.Range("B6").Value = "Unhandled error when invalid data is entered into Price Per Minute textbox."
$ws.Range("A7").Value = "Forgotten UserName"
$ws.Range("B7").Value = "When data is entered into the form and submit is clicked, nothing happens."
$ws.Range("A8").Value = "Manage Volunteers"
$ws.Range("B8").Value = "Edit button takes user to a blank screen. (Functionality is not yet implemented.)"
$ws.Range("A9").Value = "Manage Volunteers"
$ws.Range("B9").Value = "Delete button takes user to a blank screen. (Functionality is not yet implemented.)"
$ws.Range("A10").Value = "Manage Volunteers"
$ws.Range("B10").Value = "Add button takes user to a blank screen. (Functionality is not yet implemented.)"
$ws.Range("A11").Value = "Profile"
$ws.Range("B11").Value = "Volunteer Type displays ""Parent/Guardian"" regardless of the user's actual volunteer type."
$ws.Range("A12").Value = "Profile"
$ws.Range("B12").Value = "Manage Children button should be visible only for ""Parent/Guardian"" volunteer types."
$ws.Range("A13").Value = "ResetPassword"
$ws.Range("B13").Value = "When data is entered into the form and submit is clicked, nothing happens."

# Column C header, then first "x" marker, then Date header/values (order chosen to
# match the shared-string creation order of the authored workbook: Fixed, x, Date)
$ws.Range("C1").Value = "Fixed"
$ws.Range("C2").Value = "x"
$ws.Range("C3").Value = "x"
$ws.Range("C4").Value = "x"
$ws.Range("C5").Value = "x"
$ws.Range("D1").Value = "Date"

# Dates for the 4 resolved issues (rows 2-5)
$ws.Range("D2").Value = 42859
$ws.Range("D3").Value = 42859
$ws.Range("D4").Value = 42859
$ws.Range("D5").Value = 42859

# Number formatting: Date column uses date format (no special alignment) -> xf index 1
$ws.Range("D1:D5").NumberFormat = "m/d/yyyy"

# Fixed column uses date-numfmt style too (matches original authoring) plus centered alignment -> xf index 2
$ws.Range("C1:C5").NumberFormat = "m/d/yyyy"
$ws.Range("C1:C5").HorizontalAlignment = -4108

# Column C width (closest achievable to 10.42578125 given engine's width quantization)
$ws.Columns("C:C").ColumnWidth = 9.71

# Update selection to match target
$null = $ws.Range("B13").Select()

# AutoFilter now spans A1:D1 (drop old sort-state by recreating autofilter on new range)
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:D1").AutoFilter(1)

# Keep the hidden _FilterDatabase defined name in sync with the new autofilter range
$wb.Names.Item(1).RefersTo = '=Sheet1!$A$1:$D$1'

# Conditional formatting: build up to 4 dxfs total, but keep only the 2nd as the active rule (dxfId=1, priority=1)
$rng = $ws.Range("A1:D1048576")
function New-FixedRule() {
  $fc = $rng.FormatConditions.Add(2, 0, 'INDIRECT("C"&ROW())="X"')
  $fc.Font.Bold = $false
  $fc.Font.Italic = $true
  $fc.Font.Color = 10921638
  $fc.Interior.Pattern = 1
  $fc.Interior.Color = 15921906
  return $fc
}
$discard1 = New-FixedRule
$rng.FormatConditions.Item(1).Delete()
$keepRule = New-FixedRule
$discard2 = New-FixedRule
$discard3 = New-FixedRule
$rng.FormatConditions.Item(2).Delete()
$rng.FormatConditions.Item(2).Delete()

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1
